$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add "y" markers for Barometer (row 3): 3V3 (C) and GND (E)
$ws.Range("C3").Value = "y"
$ws.Range("E3").Value = "y"

# Add "y" markers for IMU (row 4): 3V3 (C), GND (E), SDA (R), SCL (S)
$ws.Range("C4").Value = "y"
$ws.Range("E4").Value = "y"
$ws.Range("R4").Value = "y"
$ws.Range("S4").Value = "y"

# Update the frozen pane / view position to match new selection
$ws.Activate()
$ws.Range("G7").Select()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("B2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("G7").Select()

$wb.Save()
